$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @{
    2 = @(44818, 200, 11000, 12000, 11500, 5750)
    3 = @(44875, 400, 7000, 7500, 7250, 3625)
    4 = @(45203, 1000, 9000, 10000, 9500, 4750)
    5 = @(44482, 240, 10000, 11000, 10500, 5250)
    6 = @(44497, 500, 9000, 10000, 9500, 4750)
    7 = @(44874, 300, 7500, 8000, 7750, 3875)
    8 = @(45211, 200, 10000, 11000, 10500, 5250)
    9 = @(44882, 440, 6000, 7000, 6500, 3250)
    10 = @(45204, 400, 9000, 10000, 9500, 4750)
    11 = @(44490, 400, 9500, 10000, 9750, 4875)
    12 = @(44517, 400, 5500, 6000, 5750, 2875)
    13 = @(44489, 160, 9500, 10000, 9750, 4875)
    14 = @(44895, 240, 3000, 3500, 3250, 1625)
    15 = @(44881, 440, 6000, 7000, 6500, 3250)
    16 = @(44475, 240, 11000, 12000, 11500, 5750)
    17 = @(44454, 160, 12000, 13000, 12500, 6250)
    18 = @(44819, 240, 11000, 12000, 11500, 5750)
    19 = @(44889, 460, 3500, 4000, 3750, 1875)
    20 = @(44455, 200, 12000, 13000, 12500, 6250)
    21 = @(44461, 200, 11000, 12000, 11500, 5750)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value = $vals[0]   # D: Fecha
    $ws.Cells.Item($row, 13).Value = $vals[1]  # M: Volumen
    $ws.Cells.Item($row, 14).Value = $vals[2]  # N: Precio minimo
    $ws.Cells.Item($row, 15).Value = $vals[3]  # O: Precio maximo
    $ws.Cells.Item($row, 16).Value = $vals[4]  # P: Precio promedio ponderado
    $ws.Cells.Item($row, 19).Value = $vals[5]  # S: Precio $/Kg
}
